# Update the "取得日時" (acquisition timestamp) column (A) for all data
# rows on the active sheet from the old timestamp to the new one.
#
# Source commit message: "Append: 2025-10-08 18:33 JST"
# The diff shows every existing row (2-18) had its column A value
# updated from 2025-10-08 18:25:50 to 2025-10-08 18:33:29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "2025-10-08 18:25:50"
$newValue = "2025-10-08 18:33:29"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
